# Apply scheduled runner updates to Behemoth_Profits sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Sheets("ALC")
$ws.Range("H28").Value = 946.5833
$ws.Range("J28").Value = 1449.5
$ws.Range("L28").Value = 1449.5
$ws.Range("N28").Value = -2419.5
$ws.Range("H68").Value = 68914.664
$ws.Range("J68").Value = 68914.664
$ws.Range("L68").Value = 68914.664
$ws.Range("N68").Value = -70412.664
$ws.Range("H71").Value = 68914.664
$ws.Range("J71").Value = 68914.664
$ws.Range("L71").Value = 206743.992
$ws.Range("N71").Value = -214231.992
$ws.Range("H80").Value = 500
$ws.Range("I80").Value = 500
$ws.Range("K80").Value = 1500
$ws.Range("M80").Value = -502
$ws.Range("H83").Value = 500
$ws.Range("I83").Value = 500
$ws.Range("K83").Value = 4500
$ws.Range("M83").Value = 492
$ws.Range("H106").Value = 2823.875
$ws.Range("I106").Value = 2427.2856
$ws.Range("J106").Value = 5600
$ws.Range("K106").Value = 2427.2856
$ws.Range("L106").Value = 5600
$ws.Range("M106").Value = -1796.2856
$ws.Range("N106").Value = -6862
$ws.Range("H137").Value = 5745.25
$ws.Range("I137").Value = 1885.2858
$ws.Range("J137").Value = 11149.2
$ws.Range("K137").Value = 5655.857400000001
$ws.Range("L137").Value = 33447.60000000001
$ws.Range("M137").Value = -3105.857400000001
$ws.Range("N137").Value = -38547.60000000001
$ws.Range("H138").Value = 2958.851
$ws.Range("J138").Value = 3144.442
$ws.Range("L138").Value = 9433.326000000001
$ws.Range("N138").Value = -19713.326

$ws = $wb.Sheets("ARM")
$ws.Range("H43").Value = 31554.334
$ws.Range("I43").Value = 27998
$ws.Range("J43").Value = 33332.5
$ws.Range("K43").Value = 27998
$ws.Range("L43").Value = 33332.5
$ws.Range("M43").Value = -27685
$ws.Range("N43").Value = -33958.5
$ws.Range("H61").Value = 25058060
$ws.Range("I61").Value = 55559150
$ws.Range("K61").Value = 55559150
$ws.Range("M61").Value = -55558938
$ws.Range("H74").Value = 14716917
$ws.Range("I74").Value = 27778460
$ws.Range("J74").Value = 22682.25
$ws.Range("K74").Value = 27778460
$ws.Range("L74").Value = 22682.25
$ws.Range("M74").Value = -27777586
$ws.Range("N74").Value = -24430.25
$ws.Range("H77").Value = 14716917
$ws.Range("I77").Value = 27778460
$ws.Range("J77").Value = 22682.25
$ws.Range("K77").Value = 138892300
$ws.Range("L77").Value = 113411.25
$ws.Range("M77").Value = -138887932
$ws.Range("N77").Value = -122147.25
$ws.Range("H101").Value = 62748.332
$ws.Range("J101").Value = 62748.332
$ws.Range("L101").Value = 62748.332
$ws.Range("N101").Value = -69238.33199999999
$ws.Range("H136").Value = 25058060
$ws.Range("I136").Value = 55559150
$ws.Range("K136").Value = 166677450
$ws.Range("M136").Value = -166674900

$ws = $wb.Sheets("BSM")
$ws.Range("H20").Value = 3890.2778
$ws.Range("I20").Value = 4388.4
$ws.Range("K20").Value = 4388.4
$ws.Range("M20").Value = -4141.4
$ws.Range("H99").Value = 2376.8235
$ws.Range("I99").Value = 2040.4
$ws.Range("K99").Value = 2040.4
$ws.Range("M99").Value = -542.4000000000001
$ws.Range("H107").Value = 2333
$ws.Range("I107").Value = 2333
$ws.Range("K107").Value = 2333
$ws.Range("M107").Value = -413
$ws.Range("H134").Value = 78072.21000000001
$ws.Range("I134").Value = 2199.4
$ws.Range("J134").Value = 120223.78
$ws.Range("K134").Value = 6598.200000000001
$ws.Range("L134").Value = 360671.34
$ws.Range("M134").Value = -4063.200000000001
$ws.Range("N134").Value = -365741.34
$ws.Range("H135").Value = 59500
$ws.Range("J135").Value = 59500
$ws.Range("L135").Value = 59500
$ws.Range("N135").Value = -69640

$ws = $wb.Sheets("CRP")
$ws.Range("H25").Value = 15000
$ws.Range("I25").Value = 15000
$ws.Range("K25").Value = 15000
$ws.Range("M25").Value = -14826
$ws.Range("H86").Value = 2999.6667
$ws.Range("I86").Value = 3000
$ws.Range("K86").Value = 3000
$ws.Range("M86").Value = -1877
$ws.Range("H89").Value = 2999.6667
$ws.Range("I89").Value = 3000
$ws.Range("K89").Value = 15000
$ws.Range("M89").Value = -9384
$ws.Range("H107").Value = 1186.5
$ws.Range("I107").Value = 1123.9231
$ws.Range("J107").Value = 2000
$ws.Range("K107").Value = 1123.9231
$ws.Range("L107").Value = 2000
$ws.Range("M107").Value = 796.0769
$ws.Range("N107").Value = -5840
$ws.Range("H134").Value = 531221.4
$ws.Range("I134").Value = 716156.9
$ws.Range("K134").Value = 2148470.7
$ws.Range("M134").Value = -2145935.7

$ws = $wb.Sheets("CUL")
$ws.Range("H7").Value = 99
$ws.Range("I7").Value = 99
$ws.Range("K7").Value = 297
$ws.Range("M7").Value = -185
$ws.Range("H20").Value = 1071.1428
$ws.Range("I20").Value = 499.8
$ws.Range("J20").Value = 2499.5
$ws.Range("K20").Value = 1499.4
$ws.Range("L20").Value = 7498.5
$ws.Range("M20").Value = -1272.4
$ws.Range("N20").Value = -7952.5
$ws.Range("H37").Value = 80999.664
$ws.Range("J37").Value = 80999.664
$ws.Range("L37").Value = 242998.992
$ws.Range("N37").Value = -243222.992
$ws.Range("H64").Value = 4200
$ws.Range("J64").Value = 4200
$ws.Range("L64").Value = 12600
$ws.Range("N64").Value = -13140
$ws.Range("H67").Value = 4200
$ws.Range("J67").Value = 4200
$ws.Range("L67").Value = 12600
$ws.Range("N67").Value = -14472
$ws.Range("H70").Value = 6640
$ws.Range("I70").Value = 4400
$ws.Range("K70").Value = 13200
$ws.Range("M70").Value = -12885
$ws.Range("H73").Value = 6640
$ws.Range("I73").Value = 4400
$ws.Range("K73").Value = 13200
$ws.Range("M73").Value = -12108
$ws.Range("H80").Value = 4454.5
$ws.Range("I80").Value = 5998
$ws.Range("J80").Value = 4283
$ws.Range("K80").Value = 17994
$ws.Range("L80").Value = 12849
$ws.Range("M80").Value = -17058
$ws.Range("N80").Value = -14721
$ws.Range("H83").Value = 4454.5
$ws.Range("I83").Value = 5998
$ws.Range("J83").Value = 4283
$ws.Range("K83").Value = 53982
$ws.Range("L83").Value = 38547
$ws.Range("M83").Value = -49302
$ws.Range("N83").Value = -47907

$ws = $wb.Sheets("GSM")
$ws.Range("H97").Value = 1609.2941
$ws.Range("I97").Value = 1204.6428
$ws.Range("K97").Value = 1204.6428
$ws.Range("M97").Value = -708.6428000000001
$ws.Range("H104").Value = 99995
$ws.Range("I104").Value = 0
$ws.Range("J104").Value = 99995
$ws.Range("K104").Value = 0
$ws.Range("L104").Value = 99995
$ws.Range("M104").ClearContents()
$ws.Range("N104").Value = -106983
$ws.Range("H122").Value = 1572.2
$ws.Range("I122").Value = 1513.2609
$ws.Range("K122").Value = 4539.7827
$ws.Range("M122").Value = -2089.7827

$ws = $wb.Sheets("LTW")
$ws.Range("H7").Value = 11827987
$ws.Range("J7").Value = 148501.72
$ws.Range("L7").Value = 148501.72
$ws.Range("N7").Value = -148725.72
$ws.Range("H40").Value = 2653.1365
$ws.Range("I40").Value = 1898.25
$ws.Range("J40").Value = 4666.1665
$ws.Range("K40").Value = 1898.25
$ws.Range("L40").Value = 4666.1665
$ws.Range("M40").Value = -1762.25
$ws.Range("N40").Value = -4938.1665
$ws.Range("H43").Value = 1087592.2
$ws.Range("I43").Value = 1842088.4
$ws.Range("J43").Value = 50160.125
$ws.Range("K43").Value = 1842088.4
$ws.Range("L43").Value = 50160.125
$ws.Range("M43").Value = -1841895.4
$ws.Range("N43").Value = -50546.125
$ws.Range("H61").Value = 999.5
$ws.Range("I61").Value = 999.5
$ws.Range("K61").Value = 999.5
$ws.Range("M61").Value = -797.5
$ws.Range("H68").Value = 2666.5
$ws.Range("I68").Value = 2499.875
$ws.Range("K68").Value = 2499.875
$ws.Range("M68").Value = -1750.875
$ws.Range("H71").Value = 2666.5
$ws.Range("I71").Value = 2499.875
$ws.Range("K71").Value = 12499.375
$ws.Range("M71").Value = -8755.375
$ws.Range("H113").Value = 999.5
$ws.Range("I113").Value = 999.5
$ws.Range("K113").Value = 999.5
$ws.Range("M113").Value = 1170.5
$ws.Range("H122").Value = 5135.8887
$ws.Range("I122").Value = 4280.409
$ws.Range("K122").Value = 12841.227
$ws.Range("M122").Value = -10391.227
$ws.Range("H126").Value = 11827987
$ws.Range("J126").Value = 148501.72
$ws.Range("L126").Value = 445505.16
$ws.Range("N126").Value = -450445.16

$ws = $wb.Sheets("WVR")
$ws.Range("H107").Value = 18519980
$ws.Range("I107").Value = 23811052
$ws.Range("J107").Value = 1229.3334
$ws.Range("K107").Value = 71433156
$ws.Range("L107").Value = 3688.0002
$ws.Range("M107").Value = -71431236
$ws.Range("N107").Value = -7528.0002
